$wb = $excel.ActiveWorkbook

# --- Metadata sheet: URL + Date text updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/mobility-alert-level"
$meta.Range("B8").Value = "2025-08-20T10:40:04+01:00"

# --- Elements sheet: Binding Value Set URL + Fixed Value (same URL as B2) text updates ---
$elem = $wb.Worksheets.Item("Elements")
$elem.Range("Z6").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/ValueSet/mobility-alert-level-vs"
$elem.Range("R5").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/mobility-alert-level"
